$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final player/position/team table (rows 2-18), reflecting the reordering
# and removal of "Wendell Carter Jr." from the roster.
$data = @(
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Kyle Kuzma", "PF", "Washington Wizards")
)

# The roster now has one fewer row (table shrinks from 18 to 17 rows),
# so remove the last row of the old range before rewriting.
$ws.Rows("19").Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
